$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.9965793528505477
$ws.Range("B5").Value = 0.9964873949579883
$ws.Range("B6").Value = 0.997342960288812
$ws.Range("B7").Value = 0.9979183673469387
$ws.Range("B8").Value = 0.9991652173913045
$ws.Range("B9").Value = 0.9991652173913045
$ws.Range("B10").Value = 0.5072535211267616
$ws.Range("B11").Value = 0.5072535211267616
$ws.Range("B12").Value = 0.8708958333333322
$ws.Range("B13").Value = 0.9943321976149965
$ws.Range("B14").Value = 0.5103176795580103
